$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "Low" / "Medium" / "High" blocks each get a small 2x2 summary
# table (header row + formula row) added in columns D/F, directly under
# their existing Min/Max/Q1/Median/Q3/IQR/Std summary blocks:
#   headerRow:  D<h> = "Mean increase"   F<h> = "Median increase"
#   dataRow:    D<h+1> = mean increase%  F<h+1> = median increase%
# plus a final block (rows 113/114) that averages the three D/F values.

$blocks = @(
    @{ Header = 18;  AvgRef = "E3";  MedRef = "E10" },
    @{ Header = 50;  AvgRef = "E35"; MedRef = "E42" },
    @{ Header = 82;  AvgRef = "E67"; MedRef = "E74" }
)

foreach ($blk in $blocks) {
    $h = $blk.Header
    $d = $h + 1

    $hdrD = $ws.Range("D$h")
    $hdrD.Value2 = "Mean increase"
    $hdrD.Font.Bold = $true

    $hdrF = $ws.Range("F$h")
    $hdrF.Value2 = "Median increase"
    $hdrF.Font.Bold = $true

    $cellD = $ws.Range("D$d")
    $cellD.Formula = "=((" + $blk.AvgRef + "/114.202998)*100)-100"
    $cellD.ClearFormats()

    $cellF = $ws.Range("F$d")
    $cellF.Formula = "=((" + $blk.MedRef + "/113.658804)*100)-100"
    $cellF.ClearFormats()
}

# Summary block: header on row 113, averages of the three blocks on row 114.
$hdrD113 = $ws.Range("D113")
$hdrD113.Value2 = "Mean increase"
$hdrD113.Font.Bold = $true

$hdrF113 = $ws.Range("F113")
$hdrF113.Value2 = "Median increase"
$hdrF113.Font.Bold = $true

$cellD114 = $ws.Range("D114")
$cellD114.Formula = "=(D19+D51+D83)/3"
$cellD114.ClearFormats()

$cellF114 = $ws.Range("F114")
$cellF114.Formula = "=(F19+F51+F83)/3"
$cellF114.ClearFormats()

# Update the active selection to match the author's final cursor position.
$ws.Range("E98").Select()
